# "bisa delte dan update password" - mark "ubah password" as done, add a new
# "middleware -> enkripsi" feature row, flag several previously-unmarked
# features as done, and add a "manage uploaded file" feature row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "ubah password" (authentication) is now done
$ws.Range("C5").Value = 1

# "middleware" row now carries an explicit not-done status
$ws.Range("C6").Value = 0

# Make room for the new "file" features (manage uploaded file) by pushing
# row 10 and below down by one row - this re-creates the gap that used to
# sit at row 7 and now sits at row 10 (between "cari file" and "file").
$ws.Rows.Item(10).Insert()

# New last feature row: "manage uploaded file" (done) - written before
# "enkripsi" below so the shared-string table keeps the same ordering as
# the authored workbook.
$ws.Range("B14").Value = "manage uploaded file"
$ws.Range("C14").Value = 1

# New middleware feature: "enkripsi" (not done yet)
$ws.Range("B7").Value = "enkripsi"
$ws.Range("C7").Value = 0

# Existing features that are now marked done
$ws.Range("C9").Value = 1
$ws.Range("C11").Value = 1
$ws.Range("C12").Value = 1
$ws.Range("C13").Value = 1

# Column B needs to widen to fit the new, longer feature names
$ws.Columns.Item(2).ColumnWidth = 20.5703125

# Move the active selection to C8, matching the author's last position
$ws.Range("C8").Select()
